$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 26.144619
$ws.Range("H2").Value = 78.433857
$ws.Range("I2").Value = 0.5211737020083955
$ws.Range("J2").Value = 0.5211737020083955
$ws.Range("M2").Value = 5.575746
$ws.Range("N2").Value = 16.727238
$ws.Range("O2").Value = 0.069238947264747
$ws.Range("P2").Value = 0.069238947264747
$ws.Range("Q2").Value = 145.775754810774
$ws.Range("R2").Value = 1311.981793296966
$ws.Range("S2").Value = 0.03608551846913226
$ws.Range("T2").Value = 0.03608551846913226
$ws.Range("G3").Value = 26.144619
$ws.Range("H3").Value = 78.433857
$ws.Range("I3").Value = 0.5211737020083955
$ws.Range("J3").Value = 0.5211737020083955
$ws.Range("O3").Value = 0.8150593598279631
$ws.Range("P3").Value = 0.815059359827963
$ws.Range("Q3").Value = 1716.026861878685
$ws.Range("R3").Value = 15444.24175690817
$ws.Range("S3").Value = 0.4247875039181325
$ws.Range("T3").Value = 0.4247875039181324
$ws.Range("G4").Value = 26.144619
$ws.Range("H4").Value = 78.433857
$ws.Range("I4").Value = 0.5211737020083955
$ws.Range("J4").Value = 0.5211737020083955
$ws.Range("M4").Value = 9.317346333333333
$ws.Range("N4").Value = 27.952039
$ws.Range("O4").Value = 0.11570169290729
$ws.Range("P4").Value = 0.11570169290729
$ws.Range("Q4").Value = 243.598469976047
$ws.Range("R4").Value = 2192.386229784423
$ws.Range("S4").Value = 0.06030067962113084
$ws.Range("T4").Value = 0.06030067962113084
$ws.Range("I5").Value = 0.3571392594830743
$ws.Range("J5").Value = 0.3571392594830742
$ws.Range("M5").Value = 5.575746
$ws.Range("N5").Value = 16.727238
$ws.Range("O5").Value = 0.069238947264747
$ws.Range("P5").Value = 0.069238947264747
$ws.Range("Q5").Value = 99.8942289741
$ws.Range("R5").Value = 899.0480607669001
$ws.Range("S5").Value = 0.02472794635351937
$ws.Range("T5").Value = 0.02472794635351937
$ws.Range("I6").Value = 0.3571392594830743
$ws.Range("J6").Value = 0.3571392594830742
$ws.Range("O6").Value = 0.8150593598279631
$ws.Range("P6").Value = 0.815059359827963
$ws.Range("S6").Value = 0.2910896962037073
$ws.Range("T6").Value = 0.2910896962037072
$ws.Range("I7").Value = 0.3571392594830743
$ws.Range("J7").Value = 0.3571392594830742
$ws.Range("M7").Value = 9.317346333333333
$ws.Range("N7").Value = 27.952039
$ws.Range("O7").Value = 0.11570169290729
$ws.Range("P7").Value = 0.11570169290729
$ws.Range("Q7").Value = 166.92817930605
$ws.Range("R7").Value = 1502.35361375445
$ws.Range("S7").Value = 0.04132161692584761
$ws.Range("T7").Value = 0.04132161692584761
$ws.Range("G8").Value = 6.104416333333333
$ws.Range("H8").Value = 18.313249
$ws.Range("I8").Value = 0.1216870385085301
$ws.Range("J8").Value = 0.1216870385085301
$ws.Range("M8").Value = 5.575746
$ws.Range("N8").Value = 16.727238
$ws.Range("O8").Value = 0.069238947264747
$ws.Range("P8").Value = 0.069238947264747
$ws.Range("Q8").Value = 34.036674952918
$ws.Range("R8").Value = 306.330074576262
$ws.Range("S8").Value = 0.008425482442095357
$ws.Range("T8").Value = 0.008425482442095357
$ws.Range("G9").Value = 6.104416333333333
$ws.Range("H9").Value = 18.313249
$ws.Range("I9").Value = 0.1216870385085301
$ws.Range("J9").Value = 0.1216870385085301
$ws.Range("O9").Value = 0.8150593598279631
$ws.Range("P9").Value = 0.815059359827963
$ws.Range("Q9").Value = 400.6691550598227
$ws.Range("R9").Value = 3606.022395538405
$ws.Range("S9").Value = 0.09918215970612328
$ws.Range("T9").Value = 0.09918215970612326
$ws.Range("G10").Value = 6.104416333333333
$ws.Range("H10").Value = 18.313249
$ws.Range("I10").Value = 0.1216870385085301
$ws.Range("J10").Value = 0.1216870385085301
$ws.Range("M10").Value = 9.317346333333333
$ws.Range("N10").Value = 27.952039
$ws.Range("O10").Value = 0.11570169290729
$ws.Range("P10").Value = 0.11570169290729
$ws.Range("Q10").Value = 56.87696114052344
$ws.Range("R10").Value = 511.8926502647109
$ws.Range("S10").Value = 0.01407939636031153
$ws.Range("T10").Value = 0.01407939636031152
